# PlanningTemplate.xlsx update
# Fills in the "Native App" planning rows (Mockups/Site/Native App/Logo/
# Moodboard/UI/DB schema/Website draft/Code) with owner (Nabil/Osamah) and
# deadline dates, matching the author's commit ("Native aangemaakt; Tabs
# mosque, prayer en favorite aangemaakt, ... ").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task data -----------------------------------------------------------
# Assign cell values in the exact order the strings were first typed so the
# shared-string table comes out in the same order as the source workbook:
# Mockups, Nabil, Site, Native App, Logo Maken, Moodboard, UI , DB schema,
# Website Draft, Osamah, Code schrijven, Website.
$ws.Range("A3").Value  = "Mockups"
$ws.Range("C3").Value  = "Nabil"
$ws.Range("B3").Value  = "Site"
$ws.Range("B4").Value  = "Native App"
$ws.Range("A5").Value  = "Logo Maken"
$ws.Range("A6").Value  = "Moodboard"
$ws.Range("A7").Value  = "UI "
$ws.Range("A8").Value  = "DB schema"
$ws.Range("A9").Value  = "Website Draft"
$ws.Range("C8").Value  = "Osamah"
$ws.Range("A10").Value = "Code schrijven"
$ws.Range("B10").Value = "Website"

# Remaining (repeat) values - string already exists so table order is kept.
$ws.Range("A4").Value  = "Mockups"
$ws.Range("C4").Value  = "Nabil"
$ws.Range("C5").Value  = "Nabil"
$ws.Range("C6").Value  = "Nabil"
$ws.Range("B7").Value  = "Native App"
$ws.Range("C7").Value  = "Nabil"
$ws.Range("C9").Value  = "Osamah"
$ws.Range("C10").Value = "Osamah"

# --- Deadline dates --------------------------------------------------------
# Apply the date number format once, then fan it out to D4:D10 via
# copy/paste-formats so every cell shares the SAME style index (matches the
# single new cellXfs entry added upstream) instead of minting a fresh style
# per cell.
$ws.Range("D3").Value = 41261
$ws.Range("D3").NumberFormat = "mm-dd-yy"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4:D10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D4").Value  = 41261
$ws.Range("D5").Value  = 41264
$ws.Range("D6").Value  = 41264
$ws.Range("D7").Value  = 41273
$ws.Range("D8").Value  = 41261
$ws.Range("D9").Value  = 41261
$ws.Range("D10").Value = 41264

# --- Selection / view -------------------------------------------------------
# The saved workbook's cursor ends up on D8 with the view scrolled back to
# the top (no frozen topLeftCell override).
$ws.Range("D8").Select() | Out-Null
